$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Swap the presentation's applied theme colours from the
#    "Integral" palette back to the standard "Office Theme" palette
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in that order).
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72

# ------------------------------------------------------------------
# 2) Re-style the table on slide 6 with the new table style id.
# ------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{1C278B7A-D8AB-42AB-9748-8E6FFA56120B}")
